$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.008085012435913
$ws.Range("B1").Value = 2.093823194503784
$ws.Range("C1").Value = 2.34978723526001
$ws.Range("D1").Value = 3.091463327407837
$ws.Range("E1").Value = 2.505070686340332
